$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B row 6: "city" -> "city_id"
$ws.Range("B6").Value = "city_id"

# New row 8: column B gets "endereco"
$ws.Range("B8").Value = "endereco"

# Update the active selection to match (cosmetic, matches diff's sheetView selection)
$ws.Range("B8").Select()
